$wb = $excel.ActiveWorkbook

# Sheet 1 = "2025"
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 973.9537847600009
$ws.Range("E2").Value = 28982.37596598056
$ws.Range("I2").Value = 16175.28135478
$ws.Range("L2").Value = 48524.529503538
$ws.Range("M2").Value = 10590.587968015
$ws.Range("N2").Value = 7169.226093134131
$ws.Range("O2").Value = 6984.121280850347

# Sheet 2 = "2030"
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 5712.560177842886
$ws.Range("E2").Value = 56106.05588781912
$ws.Range("I2").Value = 44217.8984721661
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 21984.28023276101
$ws.Range("N2").Value = 10615.03632605705
$ws.Range("O2").Value = 12072.05326959172

# Sheet 3 = "2035"
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 2861.961401238371
$ws.Range("B2").Value = 8026.889663087295
$ws.Range("E2").Value = 67297.73995507321
$ws.Range("I2").Value = 59256.42575923612
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 25464.6214365565
$ws.Range("N2").Value = 15161.1375241418
$ws.Range("O2").Value = 14770.81484578486

# Sheet 4 = "2040"
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = 2861.961401238371
$ws.Range("B2").Value = 8026.889663087295
$ws.Range("E2").Value = 67297.73995507321
$ws.Range("I2").Value = 59256.42575923612
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 25464.6214365565
$ws.Range("N2").Value = 15266.09369184379
$ws.Range("O2").Value = 14770.81484578486

# Sheet 5 = "2045"
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = 6302.873118834019
$ws.Range("B2").Value = 8026.889663087295
$ws.Range("E2").Value = 67297.73995507321
$ws.Range("I2").Value = 59256.42575923612
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 25464.6214365565
$ws.Range("N2").Value = 15804.8035822404
$ws.Range("O2").Value = 17113.37003595566

# Sheet 6 = "2050"
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = 6302.873118834019
$ws.Range("B2").Value = 8026.889663087295
$ws.Range("E2").Value = 67297.73995507321
$ws.Range("I2").Value = 59256.42575923612
$ws.Range("L2").Value = 66966.57749858923
$ws.Range("M2").Value = 25464.6214365565
$ws.Range("N2").Value = 15804.8035822404
$ws.Range("O2").Value = 17113.37003595566
